# Automatische test-sync: 2025-06-24 21:50:50
# Append the new "product Y" inquiry row to the Logs sheet and refresh the
# Dashboard category summary to reflect the updated counts.

$wb = $excel.ActiveWorkbook
$logs = $wb.Worksheets.Item("Logs")
$dash = $wb.Worksheets.Item("Dashboard")

# --- Logs sheet: append row 38 -------------------------------------------
$newRow = $logs.UsedRange.Rows.Count + 1

$logs.Cells.Item($newRow, 1).Value = "Is product Y nog op voorraad?"
$logs.Cells.Item($newRow, 2).Value = "mailmind.test@zohomail.eu"
$logs.Cells.Item($newRow, 3).Value = "Ik wil graag weten of product Y beschikbaar is."
$logs.Cells.Item($newRow, 4).Value = "Productinformatie"
$logs.Cells.Item($newRow, 5).Value = "Beste klant,`nDank je wel voor je interesse in product Y. Helaas is product Y momenteel niet op voorraad. We verwachten een nieuwe levering binnen twee weken. Mocht je op de hoogte gehouden willen worden wanneer het product weer beschikbaar is, laat het ons dan weten.`nMet vriendelijke groet,`n[Naam van het bedrijf]"
$logs.Cells.Item($newRow, 6).Value = "2025-06-24 21:50:17"
$logs.Cells.Item($newRow, 7).Value = "Ja"

# --- Extend conditional formatting ranges to include the new row ---------
$catRange = "D2:D" + $newRow
$answeredRange = "G2:G" + $newRow

$catFcs = $logs.Range("D2:D" + ($newRow - 1)).FormatConditions
$catFcs.Item(1).ModifyAppliesToRange($logs.Range($catRange))

$answeredFcs = $logs.Range("G2:G" + ($newRow - 1)).FormatConditions
$answeredFcs.Item(1).ModifyAppliesToRange($logs.Range($answeredRange))

# --- Dashboard sheet: refresh category counts -----------------------------
# "Productinformatie" now has one more reply, which bumps it to the top of
# the (count-descending) summary table, pushing the two ties below it.
$dash.Cells.Item(6, 1).Value = "Productinformatie"
$dash.Cells.Item(6, 2).Value = 3

$dash.Cells.Item(7, 1).Value = "Offerte / Prijsaanvraag"
$dash.Cells.Item(7, 2).Value = 3

$dash.Cells.Item(8, 1).Value = "Sollicitatie / Vacature"
$dash.Cells.Item(8, 2).Value = 2
